# Auto-generated edit: updates the Top-20 hot-stock ranking table
# (columns A=财联社, B=东方财富, C=同花顺) for rows 2-21 on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "至纯科技"
$ws.Range("C2").Value = "闻泰科技"

$ws.Range("A3").Value = "至纯科技"
$ws.Range("C3").Value = "北方稀土"

$ws.Range("A4").Value = "新莱应材"
$ws.Range("B4").Value = "新莱应材"
$ws.Range("C4").Value = "山子高科"

$ws.Range("B5").Value = "闻泰科技"
$ws.Range("C5").Value = "立讯精密"

$ws.Range("A6").Value = "闻泰科技"
$ws.Range("B6").Value = "凯美特气"
$ws.Range("C6").Value = "至纯科技"

$ws.Range("A7").Value = "长城军工"
$ws.Range("B7").Value = "包钢股份"
$ws.Range("C7").Value = "安泰科技"

$ws.Range("B8").Value = "安泰科技"
$ws.Range("C8").Value = "凯美特气"

$ws.Range("A9").Value = "合锻智能"
$ws.Range("B9").Value = "奥普光电"
$ws.Range("C9").Value = "海南华铁"

$ws.Range("A10").Value = "包钢股份"
$ws.Range("B10").Value = "东方财富"
$ws.Range("C10").Value = "新莱应材"

$ws.Range("A11").Value = "安泰科技"
$ws.Range("B11").Value = "长城军工"
$ws.Range("C11").Value = "合锻智能"

$ws.Range("A12").Value = "奥普光电"
$ws.Range("B12").Value = "合锻智能"
$ws.Range("C12").Value = "上海电气"

$ws.Range("A13").Value = "深振业Ａ"
$ws.Range("B13").Value = "上海电气"
$ws.Range("C13").Value = "长城军工"

$ws.Range("A14").Value = "上海电气"
$ws.Range("B14").Value = "大众公用"
$ws.Range("C14").Value = "深科技"

$ws.Range("A15").Value = "大众公用"
$ws.Range("B15").Value = "中国核建"

$ws.Range("A16").Value = "中国稀土"
$ws.Range("B16").Value = "深振业Ａ"
$ws.Range("C16").Value = "大众公用"

$ws.Range("A17").Value = "深科技"
$ws.Range("B17").Value = "宝泰隆"
$ws.Range("C17").Value = "永鼎股份"

$ws.Range("A18").Value = "中国核建"
$ws.Range("B18").Value = "深科技"
$ws.Range("C18").Value = "上海电力"

$ws.Range("A19").Value = "利和兴"
$ws.Range("B19").Value = "山子高科"
$ws.Range("C19").Value = "白银有色"

$ws.Range("A20").Value = "通富微电"
$ws.Range("B20").Value = "白银有色"
$ws.Range("C20").Value = "中电鑫龙"

$ws.Range("A21").Value = "白银有色"
$ws.Range("B21").Value = "贵州茅台"
$ws.Range("C21").Value = "金山办公"
